$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height adjustments (header rows shrunk) ---
$ws.Rows.Item(1).RowHeight = 10.8
$ws.Rows.Item(2).RowHeight = 10.8

# --- Swap the data for row 6 (Dr. ARUMUGAM J) and row 7 (Dr. REVATHY T A) ---
# Dr. REVATHY T A now appears above Dr. ARUMUGAM J.
$cols = 10
$row6vals = @()
$row7vals = @()
for ($c = 1; $c -le $cols; $c++) {
    $row6vals += ,$ws.Cells.Item(6, $c).Value2
    $row7vals += ,$ws.Cells.Item(7, $c).Value2
}
for ($c = 1; $c -le $cols; $c++) {
    $ws.Cells.Item(6, $c).Value = $row7vals[$c-1]
    $ws.Cells.Item(7, $c).Value = $row6vals[$c-1]
}

# --- Re-apply the Name/Designation font (Times New Roman, 14pt, black) across
#     the data rows (5-10) to match the refreshed formatting applied in the
#     original edit ---
$nameDesigRange = $ws.Range("A5:B10")
$nameDesigRange.Font.Name = "Times New Roman"
$nameDesigRange.Font.Size = 14
$nameDesigRange.Font.Color = 0
$nameDesigRange.Font.Bold = $false
$nameDesigRange.Font.Italic = $false
$nameDesigRange.HorizontalAlignment = -4131
